# Generate Report for Handoff
# "b.md" moved from a completed handback into a new handoff cycle:
#  - Status changes from "Handed back: in sync with en-US" to "Ready for handoff"
#  - A new handoff artifact file name / timestamp is recorded for both locales

$wb = $excel.ActiveWorkbook

function Set-HyperlinkDisplay($ws, [string]$cellAddr, [string]$newText) {
    foreach ($hl in $ws.Hyperlinks) {
        $addr = $hl.Range.Address()
        if ($addr -eq $cellAddr) {
            $hl.TextToDisplay = $newText
        }
    }
}

# ---- Overview sheet ----
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = "Ready for handoff"
$wsOverview.Range("C3").Value = "Ready for handoff"

# ---- zh-cn sheet ----
$wsZh = $wb.Worksheets.Item("zh-cn")
$wsZh.Range("B3").Value = "Ready for handoff"
$wsZh.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
Set-HyperlinkDisplay $wsZh '$C$3' "b.63290e5768f688058c7b37413b0a5c26c308f864.zh-cn.xlf"
$wsZh.Range("D3").Value = "2016-03-09 06:41:58"

# ---- de-de sheet ----
$wsDe = $wb.Worksheets.Item("de-de")
$wsDe.Range("B3").Value = "Ready for handoff"
$wsDe.Range("C3").Value = "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
Set-HyperlinkDisplay $wsDe '$C$3' "b.63290e5768f688058c7b37413b0a5c26c308f864.de-de.xlf"
$wsDe.Range("D3").Value = "2016-03-09 06:42:02"
